$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values (matching original inlineStr/text cell type).
# For Price cells whose new value is purely numeric-looking, force the
# cell to Text format first so Excel keeps storing the literal text
# (e.g. "1.003") instead of silently converting it to a number.
$ws.Cells.Item(2, 4).Value = '30.850.24'
$ws.Cells.Item(2, 5).Value = '  -0.90%  '
$ws.Cells.Item(3, 4).Value = '1.947.72'
$ws.Cells.Item(3, 5).Value = '  -0.81%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.51%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '241.49'
$ws.Cells.Item(5, 5).Value = '  -2.43%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.002'
$ws.Cells.Item(6, 5).Value = '  +0.21%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4855'
$ws.Cells.Item(7, 5).Value = '  -0.69%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2916'
$ws.Cells.Item(8, 5).Value = '  -1.60%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06834'
$ws.Cells.Item(9, 5).Value = '  -0.08%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.43'
$ws.Cells.Item(10, 5).Value = '  +1.14%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '105.01'
$ws.Cells.Item(11, 5).Value = '  -1.61%  '
$ws.Cells.Item(12, 4).Value = '1.946.82'
$ws.Cells.Item(12, 5).Value = '  -0.85%  '
$ws.Cells.Item(13, 5).Value = '  -0.28%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.289'
$ws.Cells.Item(14, 5).Value = '  -2.74%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6938'
$ws.Cells.Item(15, 5).Value = '  -3.21%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '275.39'
$ws.Cells.Item(16, 5).Value = '  -3.02%  '
$ws.Cells.Item(17, 4).Value = '30.880.95'
$ws.Cells.Item(17, 5).Value = '  -0.54%  '
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000007682'
$ws.Cells.Item(18, 5).Value = '  -0.96%  '
$ws.Cells.Item(19, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(19, 4).Value = '2.215.37'
$ws.Cells.Item(19, 5).Value = '  -0.31%  '
$ws.Cells.Item(20, 2).Value = 'Avalanche'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.09'
$ws.Cells.Item(20, 5).Value = '  -1.67%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.002'
$ws.Cells.Item(21, 5).Value = '  +0.27%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.417'
$ws.Cells.Item(22, 5).Value = '  -3.93%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.004'
$ws.Cells.Item(23, 5).Value = '  +0.67%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.445'
$ws.Cells.Item(24, 5).Value = '  -2.96%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.673'
$ws.Cells.Item(25, 5).Value = '  -3.80%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '167.36'
$ws.Cells.Item(26, 5).Value = '  -1.16%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.46'
$ws.Cells.Item(27, 5).Value = '  -2.80%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.151'
$ws.Cells.Item(28, 5).Value = '  -2.39%  '
$ws.Cells.Item(29, 5).Value = '  -3.15%  '
$ws.Cells.Item(30, 5).Value = '  -3.58%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.552'
$ws.Cells.Item(31, 5).Value = '  -2.76%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.506'
$ws.Cells.Item(32, 5).Value = '  -6.39%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.342'
$ws.Cells.Item(33, 5).Value = '  -4.10%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04837'
$ws.Cells.Item(34, 5).Value = '  -4.62%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7431'
$ws.Cells.Item(35, 5).Value = '  -3.64%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.151'
$ws.Cells.Item(36, 5).Value = '  -1.70%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.731'
$ws.Cells.Item(37, 5).Value = '  -0.04%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01979'
$ws.Cells.Item(38, 5).Value = '  -3.56%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.676'
$ws.Cells.Item(39, 5).Value = '  -1.25%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.445'
$ws.Cells.Item(40, 5).Value = '  +0.26%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '76.46'
$ws.Cells.Item(41, 5).Value = '  +3.44%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.075'
$ws.Cells.Item(42, 5).Value = '  -2.64%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.8955'
$ws.Cells.Item(43, 5).Value = '  +1.13%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '107.72'
$ws.Cells.Item(44, 5).Value = '  -1.97%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.4388'
$ws.Cells.Item(45, 5).Value = '  -1.93%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.9993'
$ws.Cells.Item(46, 5).Value = '  +0.11%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.687'
$ws.Cells.Item(47, 5).Value = '  +2.53%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '989.60'
$ws.Cells.Item(48, 5).Value = '  -0.85%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.1239'
$ws.Cells.Item(49, 5).Value = '  -2.35%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.156'
$ws.Cells.Item(50, 5).Value = '  -2.51%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '35.51'
$ws.Cells.Item(51, 5).Value = '  -1.40%  '
